# Update "想去人数" (want-to-go count) and "最低票价" (min price) figures
# plus one refreshed cover-image URL, matching a newer site scrape
# (gh-pages data refresh @ 456a3b4).
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 10900
$ws.Range("F4").Value = 262
$ws.Range("F5").Value = 1211
$ws.Range("F6").Value = 1082
$ws.Range("F7").Value = 843
$ws.Range("F8").Value = 283
$ws.Range("F9").Value = 1160
$ws.Range("F11").Value = 151
$ws.Range("F12").Value = 885
$ws.Range("F14").Value = 1986
$ws.Range("F16").Value = 980
$ws.Range("F17").Value = 827
$ws.Range("F18").Value = 554
$ws.Range("F19").Value = 804
$ws.Range("F20").Value = 916
$ws.Range("F23").Value = 91
$ws.Range("F24").Value = 630
$ws.Range("F25").Value = 648
$ws.Range("F27").Value = 353
$ws.Range("F28").Value = 1017
$ws.Range("F29").Value = 46
$ws.Range("F30").Value = 499
$ws.Range("F31").Value = 173
$ws.Range("F32").Value = 252
$ws.Range("F33").Value = 237
$ws.Range("F34").Value = 575
$ws.Range("F35").Value = 1836
$ws.Range("F37").Value = 28
$ws.Range("F38").Value = 1437
$ws.Range("F39").Value = 409
$ws.Range("F41").Value = 49
$ws.Range("F42").Value = 87
$ws.Range("F43").Value = 46
$ws.Range("F44").Value = 5
$ws.Range("F46").Value = 83
$ws.Range("F47").Value = 45
$ws.Range("F49").Value = 82

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 196
$ws.Range("G5").Value = 55
$ws.Range("F7").Value = 73
$ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202407/TxCZBf4D1721810695745.png"
$ws.Range("F14").Value = 138
$ws.Range("F15").Value = 4399

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2173
$ws.Range("F3").Value = 633
$ws.Range("F4").Value = 567

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2173
$ws.Range("F3").Value = 633
$ws.Range("F5").Value = 10900
$ws.Range("F6").Value = 262
$ws.Range("F8").Value = 567
$ws.Range("F9").Value = 1083
$ws.Range("F10").Value = 196
$ws.Range("G10").Value = 55
$ws.Range("F11").Value = 1160
$ws.Range("F13").Value = 151
$ws.Range("F14").Value = 885
$ws.Range("F15").Value = 1986
$ws.Range("F17").Value = 981
$ws.Range("F18").Value = 827
$ws.Range("F19").Value = 554
$ws.Range("F20").Value = 804
$ws.Range("F21").Value = 916
$ws.Range("F25").Value = 91
$ws.Range("F26").Value = 630
$ws.Range("F29").Value = 648
$ws.Range("F31").Value = 353
$ws.Range("F32").Value = 1017
$ws.Range("F33").Value = 46
$ws.Range("F34").Value = 499
$ws.Range("F35").Value = 173
$ws.Range("F36").Value = 252
$ws.Range("F37").Value = 237
$ws.Range("F38").Value = 28
$ws.Range("F39").Value = 1438
$ws.Range("F40").Value = 409
$ws.Range("F42").Value = 49
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 46
$ws.Range("F45").Value = 5
$ws.Range("F47").Value = 45
$ws.Range("F48").Value = 82
